$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Marking" row: marks awarded per correct answer (was 3, now 5)
$ws.Range("B11").Value = 5

# "Total" row: total correct marks (was 42 = 14*3, now 70 = 14*5)
$ws.Range("B12").Value = 70

# "Total" row: score/total max text (was "39/84", now "70/140")
$ws.Range("E12").Value = "70/140"
